# Commit "Add data for 2022-04-14": the weekly carjacking extract now runs
# through 2022-04-06 instead of 2022-04-05, so the report title (the sheet
# tab name and the matching column-header cell) is bumped, and a handful of
# neighborhood/month cells pick up the newly-tallied incidents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name) to reflect the new "through" date.
$ws.Name = "Through 2022-04-06"

# Column B's header (row 1) carries the same "through" date in its label.
$ws.Range("B1").Value = "April 2022 (through April 06)"

# Newly-populated neighborhood/month counts (previously-empty cells).
$ws.Range("F2").Value = 1    # Austin
$ws.Range("N3").Value = 1    # Englewood
$ws.Range("J5").Value = 1    # Garfield Park
$ws.Range("F12").Value = 1   # Calumet Heights
$ws.Range("N26").Value = 1   # South Shore
$ws.Range("V32").Value = 1   # Loop
$ws.Range("B44").Value = 1   # Brighton Park
$ws.Range("Z44").Value = 1   # Brighton Park
$ws.Range("F81").Value = 1   # Printers Row
$ws.Range("N91").Value = 2   # West Pullman
$ws.Range("Z91").Value = 1   # West Pullman

# Existing counts that were revised upward.
$ws.Range("V4").Value = 3    # North Lawndale: 2 -> 3
$ws.Range("V54").Value = 2   # Logan Square: 1 -> 2
